$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$co = $ws.ChartObjects()
$c2obj = $co.Item(2)
$chart = $c2obj.Chart
$ser = $chart.SeriesCollection(2)
$rng = $ws.Range("B2:B513")
Write-Host ("rng address: " + $rng.Address())
try {
  $ser.Values = $rng
  Write-Host "set values ok"
} catch {
  Write-Host ("set values err: " + $_)
}
Write-Host $ser.Formula
